{"js": "// Fixed #418 Empty AQL expressions generate empty lines.\n// Remove the empty paragraph (an empty run with no text) that sits\n// between \"... template :\" and \"End of demonstration.\" \u2014 it was being\n// emitted as a blank line for an empty AQL expression result.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const para = paragraphs.items[i];\n  if (para.text === \"\") {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fixed #418 Empty AQL expressions generate empty lines.\n# Remove the empty paragraph (no run text, other than the paragraph\n# mark) that sits between \"... template :\" and \"End of demonstration.\"\n# \u2014 it was being emitted as a blank line for an empty AQL expression\n# result.\n\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text\n    # Strip the trailing paragraph mark (and cell mark, just in case)\n    # before checking whether the paragraph has any real content.\n    $stripped = $text.TrimEnd([char]13, [char]7)\n    if ($stripped.Length -eq 0) {\n        $p.Range.Delete()\n    }\n}\n"}
